$wb = $excel.ActiveWorkbook

# --- Sheet: DP_Matrix (sheet13.xml) ---
# "Unique Binary Search Trees" (row 5) is being regrouped out of DP_Matrix
# into DP_Sequence. Delete its row so the following rows (Unique Paths,
# Unique Paths II) shift up and keep their original per-row formatting.
$wsMatrix = $wb.Worksheets.Item("DP_Matrix")
$wsMatrix.Rows("5:5").Delete() | Out-Null
$wsMatrix.PageSetup.Orientation = 1
$wsMatrix.Range("A5:XFD5").Select() | Out-Null

# --- Sheet: DP_Sequence (sheet14.xml) ---
# Re-add "Unique Binary Search Trees" as the new last row, carrying over
# the "done" (yellow) highlight it had in DP_Matrix.
$wsSeq = $wb.Worksheets.Item("DP_Sequence")
$wsSeq.Range("A7").Value = "Unique Binary Search Trees"
$wsSeq.Range("B7").Interior.ColorIndex = 6
$wsSeq.Range("A7:XFD7").Select() | Out-Null
